$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto symbol price/volume figures (GitHub Actions bot run).
# Values are written with a leading apostrophe so Excel stores them as
# literal text -- matching the original inlineStr cells -- instead of
# auto-converting the numeric-looking strings/percentages into numbers.
$updates = [ordered]@{
    'D2' = '303.18'
    'E2' = '-4.48%'
    'D3' = '35.32'
    'E3' = '-2.19%'
    'D4' = '5.074'
    'E4' = '-1.82%'
    'D5' = '0.08011'
    'E5' = '-2.89%'
    'D6' = '1.941'
    'E6' = '-9.62%'
    'D8' = '7.768'
    'E8' = '-3.24%'
    'E9' = '5.63%'
    'D10' = '0.9245'
    'E10' = '-0.29%'
    'D11' = '0.1234'
    'E11' = '19.73%'
    'D12' = '0.1861'
    'E12' = '-1.77%'
    'D13' = '0.09696'
    'E13' = '2.59%'
    'D14' = '0.03627'
    'E14' = '0.63%'
    'E15' = '-0.54%'
    'D16' = '0.001392'
    'E16' = '-3.62%'
    'D17' = '0.005687'
    'E17' = '-0.53%'
    'D18' = '3.508'
    'E18' = '1.29%'
    'D19' = '0.3406'
    'E19' = '1.07%'
    'D20' = '0.1311'
    'E20' = '-0.75%'
    'D21' = '5.044'
    'E21' = '-2.55%'
    'D22' = '0.2466'
    'E22' = '12.53%'
    'D23' = '0.04527'
    'E23' = '-1.65%'
    'D24' = '0.001217'
    'E24' = '-2.50%'
    'D25' = '0.004833'
    'E25' = '2.10%'
    'D26' = '0.0001251'
    'E26' = '-0.03%'
    'D27' = '0.0003005'
    'E27' = '-33.28%'
    'D39' = '0.01927'
    'E39' = '-3.71%'
    'D40' = '0.04719'
    'E40' = '-4.70%'
    'D41' = '0.007537'
    'E41' = '-3.62%'
    'D42' = '0.009678'
    'E42' = '23.55%'
    'E43' = '-5.09%'
    'D44' = '0.002111'
    'E44' = '-1.89%'
    'D45' = '0.01006'
    'E45' = '-14.21%'
    'D46' = '0.00006254'
    'E46' = '-3.41%'
    'D47' = '0.00000000750'
    'E47' = '-0.03%'
    'E48' = '61.63%'
    'D49' = '0.001489'
    'E49' = '-21.68%'
    'D50' = '0.00002101'
    'E50' = '-0.03%'
    'D51' = '0.0002001'
    'E51' = '-0.03%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    # Drop the quote-prefix/text number-format the line above implies so
    # the cell keeps the workbook default style, same as the source file.
    $cell.Style = "Normal"
}

